$wb = $excel.ActiveWorkbook

# Sheet 1: "Confirmation Events" - remove the "Confirmation Name" event row (row 7),
# shifting the rows below it up.
$ws1 = $wb.Worksheets.Item("Confirmation Events")
$ws1.Range("A7:F7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Sheet 2: "Candidates with events" - the "Confirmation Name" event (index 5) data is
# removed, and the data that belonged to the "Upload Sponsor Covenant" event (old index 6,
# columns Z:AA) now occupies the candidate_events.5 columns (X:Y). The now-unused trailing
# candidate_events.6 columns (Z:AA) are deleted.
$ws2 = $wb.Worksheets.Item("Candidates with events")
$ws2.Range("X4").Value = 42736
$ws2.Range("X4").NumberFormat = "m/d/yyyy"
$ws2.Range("Y4").Value = $false
$ws2.Range("Z1:AA10").EntireColumn.Delete()
